$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header values to add into the previously-blank row 1 (A1:EP1)
$headers = @(
    '_id',
    'Category',
    'Topic',
    'Data Source',
    'Characteristic',
    'City of Toronto',
    'Agincourt North',
    'Agincourt South-Malvern West',
    'Alderwood',
    'Annex',
    'Banbury-Don Mills',
    'Bathurst Manor',
    'Bay Street Corridor',
    'Bayview Village',
    'Bayview Woods-Steeles',
    'Bedford Park-Nortown',
    'Beechborough-Greenbrook',
    'Bendale',
    'Birchcliffe-Cliffside',
    'Black Creek',
    'Blake-Jones',
    'Briar Hill-Belgravia',
    'Bridle Path-Sunnybrook-York Mills',
    'Broadview North',
    'Brookhaven-Amesbury',
    'Cabbagetown-South St. James Town',
    'Caledonia-Fairbank',
    'Casa Loma',
    'Centennial Scarborough',
    'Church-Yonge Corridor',
    'Clairlea-Birchmount',
    'Clanton Park',
    'Cliffcrest',
    'Corso Italia-Davenport',
    'Danforth',
    'Danforth East York',
    'Don Valley Village',
    'Dorset Park',
    'Dovercourt-Wallace Emerson-Junction',
    'Downsview-Roding-CFB',
    'Dufferin Grove',
    'East End-Danforth',
    'Edenbridge-Humber Valley',
    'Eglinton East',
    'Elms-Old Rexdale',
    'Englemount-Lawrence',
    'Eringate-Centennial-West Deane',
    'Etobicoke West Mall',
    'Flemingdon Park',
    'Forest Hill North',
    'Forest Hill South',
    'Glenfield-Jane Heights',
    'Greenwood-Coxwell',
    'Guildwood',
    'Henry Farm',
    'High Park North',
    'High Park-Swansea',
    'Highland Creek',
    'Hillcrest Village',
    'Humber Heights-Westmount',
    'Humber Summit',
    'Humbermede',
    'Humewood-Cedarvale',
    'Ionview',
    'Islington-City Centre West',
    'Junction Area',
    'Keelesdale-Eglinton West',
    'Kennedy Park',
    'Kensington-Chinatown',
    'Kingsview Village-The Westway',
    'Kingsway South',
    'Lambton Baby Point',
    'L''Amoreaux',
    'Lansing-Westgate',
    'Lawrence Park North',
    'Lawrence Park South',
    'Leaside-Bennington',
    'Little Portugal',
    'Long Branch',
    'Malvern',
    'Maple Leaf',
    'Markland Wood',
    'Milliken',
    'Mimico (includes Humber Bay Shores)',
    'Morningside',
    'Moss Park',
    'Mount Dennis',
    'Mount Olive-Silverstone-Jamestown',
    'Mount Pleasant East',
    'Mount Pleasant West',
    'New Toronto',
    'Newtonbrook East',
    'Newtonbrook West',
    'Niagara',
    'North Riverdale',
    'North St. James Town',
    'Oakridge',
    'Oakwood Village',
    'O''Connor-Parkview',
    'Old East York',
    'Palmerston-Little Italy',
    'Parkwoods-Donalda',
    'Pelmo Park-Humberlea',
    'Playter Estates-Danforth',
    'Pleasant View',
    'Princess-Rosethorn',
    'Regent Park',
    'Rexdale-Kipling',
    'Rockcliffe-Smythe',
    'Roncesvalles',
    'Rosedale-Moore Park',
    'Rouge',
    'Runnymede-Bloor West Village',
    'Rustic',
    'Scarborough Village',
    'South Parkdale',
    'South Riverdale',
    'St.Andrew-Windfields',
    'Steeles',
    'Stonegate-Queensway',
    'Tam O''Shanter-Sullivan',
    'Taylor-Massey',
    'The Beaches',
    'Thistletown-Beaumond Heights',
    'Thorncliffe Park',
    'Trinity-Bellwoods',
    'University',
    'Victoria Village',
    'Waterfront Communities-The Island',
    'West Hill',
    'West Humber-Clairville',
    'Westminster-Branson',
    'Weston',
    'Weston-Pelham Park',
    'Wexford/Maryvale',
    'Willowdale East',
    'Willowdale West',
    'Willowridge-Martingrove-Richview',
    'Woburn',
    'Woodbine Corridor',
    'Woodbine-Lumsden',
    'Wychwood',
    'Yonge-Eglinton',
    'Yonge-St.Clair',
    'York University Heights',
    'Yorkdale-Glen Park'
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Select the whole of row 1 (full-row selection, A1:XFD1) and scroll so
# that column B becomes the left-most visible column
$ws.Rows.Item(1).Select()
$ws.Application.ActiveWindow.ScrollColumn = 2
